# ---------------------------------------------------------------------------
# Applies the edits captured by the target diff to the "variables" sheet of
# data_extraction_setSD_checkedTR.xlsx:
#   - A new coded answer "exclude_NA" was added to the proxy_decision (AC)
#     column's pick-list and re-applied to the rows that previously used the
#     plain "exclude" answer for the female/adult body-condition proxies.
#   - Several repeated_trait_ID-style tally columns (I/J) and the "sex" (AO)
#     column were recounted/relabelled for a block of rows.
#   - Column H and AA were widened to fit the now-longer entries.
#   - The user's last selection moved from AD78 to AF44.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

# --- Column width adjustments (H and AA widened) ---
$ws.Columns.Item(8).ColumnWidth = 13
$ws.Columns.Item(27).ColumnWidth = 32

# --- proxy_decision (AC): "exclude" -> "exclude_NA" ---
$acRows = 14,15,16,17,18,63,64,65,66,67,68,69,70,71,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95
foreach ($r in $acRows) {
    $ws.Range("AC$r").Value = "exclude_NA"
}

# --- I/J tally columns, rows 42-62 ---
$ws.Range("I42").Value = 2
$ws.Range("J42").Value = 2
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 2
$ws.Range("J45").Value = 3
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = 4
$ws.Range("J47").Value = 3
$ws.Range("I48").Value = 2
$ws.Range("J48").Value = 4
$ws.Range("J49").Value = 5
$ws.Range("I50").Value = 2
$ws.Range("J50").Value = 6
$ws.Range("J51").Value = 5
$ws.Range("I52").Value = 2
$ws.Range("J52").Value = 6
$ws.Range("J53").Value = 7
$ws.Range("I54").Value = 2
$ws.Range("J54").Value = 8
$ws.Range("J55").Value = 7
$ws.Range("I56").Value = 2
$ws.Range("J56").Value = 8
$ws.Range("J57").Value = 9
$ws.Range("I58").Value = 2
$ws.Range("J58").Value = 10
$ws.Range("J59").Value = 9
$ws.Range("I60").Value = 2
$ws.Range("J60").Value = 10
$ws.Range("J61").Value = 11
$ws.Range("I62").Value = 2
$ws.Range("J62").Value = 12

# --- J tally column, rows 84-98 (shifted down by 3) ---
$ws.Range("J84").Value = 19
$ws.Range("J85").Value = 20
$ws.Range("J86").Value = 21
$ws.Range("J87").Value = 22
$ws.Range("J88").Value = 23
$ws.Range("J89").Value = 24
$ws.Range("J90").Value = 25
$ws.Range("J91").Value = 26
$ws.Range("J92").Value = 27
$ws.Range("J93").Value = 28
$ws.Range("J94").Value = 29
$ws.Range("J95").Value = 30
$ws.Range("J96").Value = 31
$ws.Range("J97").Value = 32
$ws.Range("J98").Value = 33

# --- J tally + AO ("sex": NA -> b) columns, rows 107-113 ---
$ws.Range("J107").Value = 2
$ws.Range("AO107").Value = "b"
$ws.Range("J108").Value = 3
$ws.Range("AO108").Value = "b"
$ws.Range("J109").Value = 4
$ws.Range("AO109").Value = "b"
$ws.Range("J110").Value = 5
$ws.Range("AO110").Value = "b"
$ws.Range("J111").Value = 6
$ws.Range("AO111").Value = "b"
$ws.Range("J112").Value = 7
$ws.Range("AO112").Value = "b"
$ws.Range("J113").Value = 8
$ws.Range("AO113").Value = "b"

# --- J tally column, rows 115-119 ---
$ws.Range("J115").Value = 2
$ws.Range("J116").Value = 3
$ws.Range("J117").Value = 4
$ws.Range("J118").Value = 5
$ws.Range("J119").Value = 6

# --- Restore the author's final selection ---
$ws.Activate()
$ws.Range("AF44").Select()
